# Regenerate save_data to use K (strikeouts) instead of Strike# (pitch-count-
# based strikeout proxy) in column G. Only the numeric values of column G
# change; the header label ("K") is already correct and stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 8
    4  = 9
    5  = 10
    6  = 6
    7  = 5
    8  = 7
    9  = 15
    10 = 6
    11 = 7
    12 = 7
    13 = 5
    14 = 9
    15 = 5
    16 = 8
    17 = 11
    18 = 12
    19 = 1
    20 = 9
    21 = 7
    22 = 2
    23 = 7
    24 = 9
    25 = 8
    26 = 4
    27 = 10
    28 = 7
    29 = 4
    30 = 10
    31 = 7
    32 = 5
    33 = 6
    34 = 5
    35 = 9
    36 = 1
    37 = 3
    38 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
